# ran rest of March TP2 run 1 samples and last sample from run5
# Append one new CRM-accuracy sample (row 36) below the existing data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample row: date, CRM value, batch value, % off (formula), batch #, notes
$ws.Range("A36").Value = 20210525
$ws.Range("B36").Value = 2223.5051632014602
$ws.Range("C36").Value = 2224.4699999999998
$ws.Range("D36").Formula = "=100*(B36-C36)/C36"
$ws.Range("E36").Value = 180
$ws.Range("F36").Value = "CRM opened 20210418"

# Reposition the view/selection the way the author left the sheet after entry.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F40").Select()
